$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("AD3").Value = 8.5
$ws.Range("AJ3").Value = 6
$ws.Range("AL3").Value = 81
$ws.Range("AN3").Value = 17

# Row 4
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("U4").Value = 5
$ws.Range("V4").Value = 1.17
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 1.11

# Row 5
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 2.4
$ws.Range("L5").Value = 3.4
$ws.Range("U5").Value = 4.8
$ws.Range("AA5").Value = 2.25
$ws.Range("AB5").Value = 1.57
$ws.Range("AC5").Value = 7
$ws.Range("AN5").Value = 10
$ws.Range("AO5").Value = 11
$ws.Range("AQ5").Value = 26

# Row 8
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.5
$ws.Range("J8").Value = 2.63
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 4
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("S8").Value = 1.83
$ws.Range("T8").Value = 2.03
$ws.Range("W8").Value = 3
$ws.Range("X8").Value = 1.36
$ws.Range("Y8").Value = 1.36
$ws.Range("Z8").Value = 3
$ws.Range("AA8").Value = 1.62
$ws.Range("AB8").Value = 2.2
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 10
$ws.Range("AF8").Value = 17
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 12
$ws.Range("AJ8").Value = 6.5
$ws.Range("AK8").Value = 13
$ws.Range("AN8").Value = 19
$ws.Range("AR8").Value = 34
$ws.Range("AS8").Value = 151

# Row 9
$ws.Range("S9").Value = 1.53
$ws.Range("T9").Value = 2.4
$ws.Range("U9").Value = 1.88
$ws.Range("V9").Value = 1.93
$ws.Range("W9").Value = 2.25
$ws.Range("X9").Value = 1.57

# Row 10
$ws.Range("G10").Value = 2.75
$ws.Range("H10").Value = 3.2
$ws.Range("K10").Value = 2.05
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75
$ws.Range("S10").Value = 2.25
$ws.Range("T10").Value = 1.62
$ws.Range("W10").Value = 4
$ws.Range("X10").Value = 1.22
$ws.Range("AA10").Value = 1.91
$ws.Range("AB10").Value = 1.91
$ws.Range("AC10").Value = 8
$ws.Range("AG10").Value = 23
$ws.Range("AH10").Value = 34
$ws.Range("AI10").Value = 8.5
$ws.Range("AM10").Value = 8
$ws.Range("AN10").Value = 12
$ws.Range("AQ10").Value = 23
$ws.Range("AR10").Value = 34
$ws.Range("AS10").Value = 301

# Row 12
$ws.Range("H12").Value = 3.25
$ws.Range("I12").Value = 2.2
$ws.Range("K12").Value = 2.05
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.4
$ws.Range("P12").Value = 2.75
$ws.Range("S12").Value = 2.25
$ws.Range("T12").Value = 1.62
$ws.Range("W12").Value = 4
$ws.Range("X12").Value = 1.22
$ws.Range("AA12").Value = 1.91
$ws.Range("AB12").Value = 1.91
$ws.Range("AC12").Value = 9.5
$ws.Range("AH12").Value = 41
$ws.Range("AI12").Value = 8.5
$ws.Range("AJ12").Value = 6
$ws.Range("AM12").Value = 7
$ws.Range("AO12").Value = 9.5
$ws.Range("AP12").Value = 21
$ws.Range("AQ12").Value = 19
$ws.Range("AS12").Value = 301

# Row 15
$ws.Range("G15").Value = 5.75
$ws.Range("H15").Value = 3.2
$ws.Range("I15").Value = 1.73
$ws.Range("J15").Value = 5.5
$ws.Range("K15").Value = 2.05
$ws.Range("L15").Value = 2.4
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 7.5
$ws.Range("S15").Value = 2.25
$ws.Range("T15").Value = 1.62
$ws.Range("W15").Value = 4
$ws.Range("X15").Value = 1.22
$ws.Range("Y15").Value = 1.5
$ws.Range("Z15").Value = 2.5
$ws.Range("AC15").Value = 12
$ws.Range("AD15").Value = 26
$ws.Range("AE15").Value = 19
$ws.Range("AI15").Value = 7.5
$ws.Range("AJ15").Value = 6.5
$ws.Range("AK15").Value = 19
$ws.Range("AN15").Value = 7
$ws.Range("AO15").Value = 9
$ws.Range("AP15").Value = 13
$ws.Range("AQ15").Value = 17

# Row 18
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 1.65
$ws.Range("J18").Value = 4.45
$ws.Range("K18").Value = 2.37
$ws.Range("L18").Value = 2.15
$ws.Range("S18").Value = 1.52
$ws.Range("T18").Value = 2.37
$ws.Range("AB18").Value = 2.3
$ws.Range("AC18").Value = 16.5
$ws.Range("AD18").Value = 28
$ws.Range("AF18").Value = 70
$ws.Range("AJ18").Value = 8.25
$ws.Range("AK18").Value = 13
$ws.Range("AL18").Value = 45
$ws.Range("AM18").Value = 10
$ws.Range("AN18").Value = 10
$ws.Range("AP18").Value = 14.5

# Row 19
$ws.Range("G19").Value = 2.75
$ws.Range("I19").Value = 2.8
$ws.Range("J19").Value = 3.5
$ws.Range("K19").Value = 1.91
$ws.Range("L19").Value = 3.6
$ws.Range("AD19").Value = 12
$ws.Range("AF19").Value = 29
$ws.Range("AM19").Value = 7
$ws.Range("AN19").Value = 12
$ws.Range("AO19").Value = 11
$ws.Range("AQ19").Value = 26

# Row 20
$ws.Range("G20").Value = 5.4
$ws.Range("H20").Value = 4.1
$ws.Range("I20").Value = 1.53
$ws.Range("J20").Value = 5.2
$ws.Range("K20").Value = 2.32
$ws.Range("L20").Value = 2.05
$ws.Range("N20").Value = 8.25
$ws.Range("O20").Value = 1.22
$ws.Range("P20").Value = 3.8
$ws.Range("S20").Value = 1.65
$ws.Range("T20").Value = 2.1
$ws.Range("W20").Value = 2.6
$ws.Range("X20").Value = 1.44
$ws.Range("AA20").Value = 1.78
$ws.Range("AB20").Value = 1.93
$ws.Range("AC20").Value = 16
$ws.Range("AD20").Value = 32
$ws.Range("AE20").Value = 17
$ws.Range("AF20").Value = 100
$ws.Range("AG20").Value = 50
$ws.Range("AI20").Value = 8.25
$ws.Range("AJ20").Value = 8
$ws.Range("AK20").Value = 16
$ws.Range("AM20").Value = 7.6
$ws.Range("AN20").Value = 7.7
$ws.Range("AP20").Value = 11
$ws.Range("AQ20").Value = 11.75
$ws.Range("AS20").Value = 500
